$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 45413
$ws.Range("B3").Value = 141125
$ws.Range("C3").Value = 13810475
$ws.Range("D3").Value = 97.859875996457035
$ws.Range("E3").Value = 1254.9000000000001
$ws.Range("F3").Value = 264873.09999999998
$ws.Range("G3").Value = 211.07108136106459
$ws.Range("K3").Value = 142379.9
$ws.Range("L3").Value = 14075348.1
$ws.Range("M3").Value = 98.857690586943804
$ws.Range("N3").Value = 969809.46
$ws.Range("O3").Value = 252347.47
$ws.Range("P3").Value = 0.26020314340922185
$ws.Range("Q3").Value = 298.25
$ws.Range("R3").Value = 56047.8
$ws.Range("S3").Value = 187.9222129086337
$ws.Range("T3").Value = 1200
$ws.Range("U3").Value = 195285.6
$ws.Range("V3").Value = 162.738
$ws.Range("W3").Value = 5918
$ws.Range("X3").Value = 759300
$ws.Range("Y3").Value = 128.30348090571138
$ws.Range("Z3").Value = 7416.25
$ws.Range("AA3").Value = 1010633.4
$ws.Range("AB3").Value = 136.27283330524187
$ws.Range("AC3").Value = 13264.519999999999
$ws.Range("AD3").Value = 3018545.8
$ws.Range("AE3").Value = 227.56540002955253
$ws.Range("AF3").Value = 11150
$ws.Range("AG3").Value = 1378750
$ws.Range("AH3").Value = 123.65470852017937
$ws.Range("AR3").Value = 7416.25
$ws.Range("AS3").Value = 1010633.4
$ws.Range("AT3").Value = 136.27283330524187

$ws.Range("A4").Select()
